$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 462; this shifts the existing rows 462..537 down to 463..538
# and copies formatting (e.g. the date style on column D) from the row being split.
$ws.Rows(462).Insert()

# Populate the newly inserted row 462 with a new data record (same shape as the
# surrounding rows: Agrícola del Norte S.A. de Arica / Arica y Parinacota / Zanahoria).
$ws.Range("A462").Value = 1
$ws.Range("B462").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C462").Value = "Arica y Parinacota"
$ws.Range("D462").Value = 45077
$ws.Range("E462").Value = 15
$ws.Range("F462").Value = 100114013
$ws.Range("G462").Value = "Zanahoria"
$ws.Range("H462").Value = "Sin especificar"
$ws.Range("I462").Value = "Primera"
$ws.Range("J462").Value = 30
$ws.Range("K462").Value = 12000
$ws.Range("L462").Value = 13000
$ws.Range("M462").Value = 12333
$ws.Range("N462").Value = "$/saco 25 kilos"
$ws.Range("O462").Value = "Región de Arica y Parinacota"
$ws.Range("P462").Value = 493
$ws.Range("Q462").Value = 25
$ws.Range("R462").Value = "Hortaliza"
